# Update existing rows 2-4 and add new rows 5-6 per the diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('A2').Value = 'Divyani Jaiswal '
$ws.Range('B2').Value = 'divyani.jaiswal@ltimindtree.com'
$ws.Range('C2').Value = 'COD'
$ws.Range('D2').Value = 0.3
$ws.Range('E2').Value = 30
$ws.Range('F2').Value = '2025-09-22 | 05:59:42 PM'
$ws.Range('I2').Value = 'The provided Java program seems to have issues with calculating the index of alphabetic characters in a given string, as evident from the logs showing comparison failures across various test cases. The program appears to be incorrectly calculating the index, often resulting in negative values or incorrect positions, indicating a logical implementation error. The expected output and actual output differences suggest that the indexing calculation is not correctly based on the standard ordering of the alphabet.'
$ws.Range('J2').Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX19DguVyvWUWUtu3EQBQd%2F0zsUNCzoSuoVpC57Ias7H%2F8kmozU1A2hi0Q8edBw7mZ1NFryNp7AZL1bMqn4Z%2FBKUkVz36%2BP5jDJVZDIE3Iky1z5SusHmcbahKGxrbNcfCVbYjvM7yO7LE%2Bg%3D%3D'

# Row 3
$ws.Range('A3').Value = 'NITISH KUMAR GUPTA'
$ws.Range('B3').Value = 'nitish.gupta2@ltimindtree.com'
$ws.Range('C3').Value = 'COD'
$ws.Range('D3').Value = 16.5
$ws.Range('E3').Value = 30
$ws.Range('F3').Value = '2025-03-25 | 05:18:57 PM'
$ws.Range('I3').Value = 'Based on the logs and description provided, it appears that there are issues with the HTTP request methods and student data retrieval/deletion. The test cases `testGetAllStudents` and `testDeleteStudentById` are failing with status code mismatches (405 and 404 respectively), indicating potential problems with the REST API endpoint configurations or the database interactions in the `StudentController` and `StudentService` classes.'
$ws.Range('J3').Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2BmmBhxr3%2Fh8dqT%2FDWerpcdT%2BJ88UuBsYbWoQdrU%2BkeRCP23RHl5rtNhBRiimvjsD%2F5BvgDQ339%2FXfHMIhWQceOmkhTde3VlSHaeapL2rJDCpc9Jo3Sxxyfv9L1wYPHY9nbkbazwuV2gw%3D%3D'

# Row 4
$ws.Range('A4').Value = 'Sreeja Reddy Minumula'
$ws.Range('B4').Value = 'minumula.sreejareddy@ltimindtree.com'
$ws.Range('C4').Value = 'COD'
$ws.Range('D4').Value = 31.200000000000003
$ws.Range('E4').Value = 60
$ws.Range('F4').Value = '2025-06-20 | 05:34:27 PM'
$ws.Range('I4').Value = 'The logs indicate that there are two test failures, `testVerifyTitleExistsInPagesPackage` and `testHandleDropdownWithSelect`, both reporting that no method named `verifyTitle` was found in any file within the Pages package. This suggests that the `verifyTitle` method is missing or not properly defined in the Pages package, causing the test cases to fail. 
However, as per description, `verifyTitle` method should be created in `utils/WebDriverHelper` class not in Pages package. Therefore, it can be assumed that the method is probably created at wrong location. 
Also, it seems like `handleDropdown` method is also not properly implemented or utilized as expected. 
Therefore, corrections are needed in terms of method implementation and its location.'
$ws.Range('J4').Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX19JaIbFVNlrvNpAOgu3ctM6j7niyokXCGUBZuZfB%2FKuckU3XS4o9o%2Fd1W5qvfJjRgsuNH0Pb2TZ5QnhXqvSTp6K0fRXSOGLAWVhP0I2YSzHUls9EE0BbXj74ft3gGteu0nRQpehOeLTKw%3D%3D'

# Row 5
$ws.Range('A5').Value = 'Akash Rai'
$ws.Range('B5').Value = 'akash.rai2@ltimindtree.com'
$ws.Range('C5').Value = 'COD1'
$ws.Range('D5').Value = 9
$ws.Range('E5').Value = 30
$ws.Range('F5').Value = '2025-09-12 | 04:55:32 PM'
$ws.Range('I5').Value = 'Based on the logs and description provided, the Turf Management System in C# appears to have several issues, including `NullReferenceException` errors when adding, deleting, and updating turfs, indicating that some objects are not being properly initialized. Additionally, the system seems to have logical errors in handling menu options, displaying turfs, and updating/deleting turfs, resulting in incorrect output or error messages. These issues suggest that the system requires debugging and refinement to ensure proper functionality.'
$ws.Range('J5').Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2BtqOsiVP9Frfbb4%2Fj%2BH1%2FSg2fBaBUBoZkdx9wqSdyHn3R1ydHTw%2Fz29frEkeu5pLe0CW%2F%2BqRMCuPXInIq0h2lM597u18T6YBlJvGrnHx7T2xyGtF2GvPPsm2uSwf1Xpf03G431NU9QgQ%3D%3D'

# Row 6
$ws.Range('A6').Value = 'Akash Rai'
$ws.Range('B6').Value = 'akash.rai2@ltimindtree.com'
$ws.Range('C6').Value = 'COD2'
$ws.Range('D6').Value = 12
$ws.Range('E6').Value = 30
$ws.Range('F6').Value = '2025-09-12 | 04:55:32 PM'
$ws.Range('I6').Value = 'Based on the logs and description provided, the issues in the Kabaddi Team Management System include data type mismatch errors, incorrect syntax in data queries, and handling of edge cases such as no records found. The logs indicate specific test failures, including inserting records with incorrect data types, incorrect display of players above a points threshold, and errors in deleting players not containing a specific word. 
 To fix these issues, the code should be reviewed for correct data type usage, query syntax, and edge case handling. 
Likewise,  here are some potential solutions:
*   **Test_AddPlayer_Should_Insert_Record**: Ensure that the `MatchesPlayed` column is being assigned an integer value, not a string like "Raider".
*   **Test_DisplayPlayersAbovePointsThreshold_Should_Output_Records**: Verify that the query to display players above a points threshold is correct and that the expected records are being returned.
*   **Test_DeletePlayersNotContainingWord_Should_Remove_Record**: Check that the syntax for the delete query is correct and that the word being searched for is not causing any issues.
*   **Test_DisplayPlayersAbovePointsThreshold_Should_Handle_No_Records_Found**: Add a check to handle the case where no records are found above the points threshold.
*   **Test_DeletePlayersNotContainingWord_Should_Handle_No_Record_Found**: Modify the delete query to handle cases where no records are found not containing the specified word.
*   **Test_UpdatePlayerDetails_Should_Handle_No_Record_Found**: Update the update query to handle cases where no records are found for the given player name.'
$ws.Range('J6').Value = 'https://admin.ltimindtree.iamneo.ai/result?testId=U2FsdGVkX1%2BtqOsiVP9Frfbb4%2Fj%2BH1%2FSg2fBaBUBoZkdx9wqSdyHn3R1ydHTw%2Fz29frEkeu5pLe0CW%2F%2BqRMCuPXInIq0h2lM597u18T6YBlJvGrnHx7T2xyGtF2GvPPsm2uSwf1Xpf03G431NU9QgQ%3D%3D'
